$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of diary entries for Murthy Routhula working on Tableau
$dates = @(44626, 44628, 44629, 44630, 44631)
$startRow = 8

$ws.Range("C3").Copy()

for ($i = 0; $i -lt $dates.Length; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = "Murthy Routhula"
    $ws.Cells.Item($r, 2).Value = "Worked on Tableau"
    $ws.Cells.Item($r, 3).Value = $dates[$i]
    $ws.Cells.Item($r, 3).PasteSpecial(-4122)
}

$ws.Range("B14").Select()
